$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H107").Value = 326.77777
$ws_ALC.Range("I107").Value = 235.9
$ws_ALC.Range("J107").Value = 586.4286
$ws_ALC.Range("K107").Value = 235.9
$ws_ALC.Range("L107").Value = 586.4286
$ws_ALC.Range("M107").Value = 1684.1
$ws_ALC.Range("N107").Value = -4426.4286
$ws_ALC.Range("H116").Value = 3280.4546
$ws_ALC.Range("I116").Value = 4041
$ws_ALC.Range("J116").Value = 2646.6667
$ws_ALC.Range("K116").Value = 4041
$ws_ALC.Range("L116").Value = 2646.6667
$ws_ALC.Range("M116").Value = -599
$ws_ALC.Range("N116").Value = -9530.6667

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H2").Value = 125861.375
$ws_ARM.Range("I2").Value = 333963.66
$ws_ARM.Range("J2").Value = 1000
$ws_ARM.Range("K2").Value = 333963.66
$ws_ARM.Range("L2").Value = 1000
$ws_ARM.Range("M2").Value = -333850.66
$ws_ARM.Range("N2").Value = -1226
$ws_ARM.Range("H32").Value = 4306.5396
$ws_ARM.Range("I32").Value = 3415
$ws_ARM.Range("J32").Value = 9518.615
$ws_ARM.Range("K32").Value = 3415
$ws_ARM.Range("L32").Value = 9518.615
$ws_ARM.Range("M32").Value = -3128
$ws_ARM.Range("N32").Value = -10092.615
$ws_ARM.Range("H74").Value = 2068.3057
$ws_ARM.Range("I74").Value = 645.8246
$ws_ARM.Range("J74").Value = 7473.7334
$ws_ARM.Range("K74").Value = 645.8246
$ws_ARM.Range("L74").Value = 7473.7334
$ws_ARM.Range("M74").Value = 228.1754
$ws_ARM.Range("N74").Value = -9221.733400000001
$ws_ARM.Range("H77").Value = 2068.3057
$ws_ARM.Range("I77").Value = 645.8246
$ws_ARM.Range("J77").Value = 7473.7334
$ws_ARM.Range("K77").Value = 3229.123
$ws_ARM.Range("L77").Value = 37368.667
$ws_ARM.Range("M77").Value = 1138.877
$ws_ARM.Range("N77").Value = -46104.667
$ws_ARM.Range("H97").Value = 50433
$ws_ARM.Range("I97").Value = 59282.824
$ws_ARM.Range("J97").Value = 284
$ws_ARM.Range("K97").Value = 59282.824
$ws_ARM.Range("L97").Value = 284
$ws_ARM.Range("M97").Value = -58786.824
$ws_ARM.Range("N97").Value = -1276
$ws_ARM.Range("H110").Value = 718.3125
$ws_ARM.Range("I110").Value = 565.2727
$ws_ARM.Range("J110").Value = 1055
$ws_ARM.Range("K110").Value = 565.2727
$ws_ARM.Range("L110").Value = 1055
$ws_ARM.Range("M110").Value = 1479.7273
$ws_ARM.Range("N110").Value = -5145
$ws_ARM.Range("H116").Value = 125861.375
$ws_ARM.Range("I116").Value = 333963.66
$ws_ARM.Range("J116").Value = 1000
$ws_ARM.Range("K116").Value = 333963.66
$ws_ARM.Range("L116").Value = 1000
$ws_ARM.Range("M116").Value = -331669.66
$ws_ARM.Range("N116").Value = -5588
$ws_ARM.Range("H132").Value = 3785.831
$ws_ARM.Range("I132").Value = 2694.8655
$ws_ARM.Range("J132").Value = 6771.6313
$ws_ARM.Range("K132").Value = 8084.5965
$ws_ARM.Range("L132").Value = 20314.8939
$ws_ARM.Range("M132").Value = -5554.5965
$ws_ARM.Range("N132").Value = -25374.8939

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H3").Value = 125861.375
$ws_BSM.Range("I3").Value = 333963.66
$ws_BSM.Range("J3").Value = 1000
$ws_BSM.Range("K3").Value = 333963.66
$ws_BSM.Range("L3").Value = 1000
$ws_BSM.Range("M3").Value = -333849.66
$ws_BSM.Range("N3").Value = -1228
$ws_BSM.Range("H134").Value = 702.4194
$ws_BSM.Range("I134").Value = 632.68604
$ws_BSM.Range("J134").Value = 1559.1428
$ws_BSM.Range("K134").Value = 1898.05812
$ws_BSM.Range("L134").Value = 4677.428400000001
$ws_BSM.Range("M134").Value = 636.9418799999999
$ws_BSM.Range("N134").Value = -9747.428400000001

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H16").Value = 1724.409
$ws_CRP.Range("I16").Value = 1075.8
$ws_CRP.Range("J16").Value = 3114.2856
$ws_CRP.Range("K16").Value = 1075.8
$ws_CRP.Range("L16").Value = 3114.2856
$ws_CRP.Range("M16").Value = -788.8
$ws_CRP.Range("N16").Value = -3688.2856
$ws_CRP.Range("H113").Value = 1724.409
$ws_CRP.Range("I113").Value = 1075.8
$ws_CRP.Range("J113").Value = 3114.2856
$ws_CRP.Range("K113").Value = 1075.8
$ws_CRP.Range("L113").Value = 3114.2856
$ws_CRP.Range("M113").Value = 1094.2
$ws_CRP.Range("N113").Value = -7454.2856
$ws_CRP.Range("H132").Value = 15154920
$ws_CRP.Range("I132").Value = 23813714
$ws_CRP.Range("J132").Value = 2030.25
$ws_CRP.Range("K132").Value = 71441142
$ws_CRP.Range("L132").Value = 6090.75
$ws_CRP.Range("M132").Value = -71438612
$ws_CRP.Range("N132").Value = -11150.75

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H130").Value = 2083.3333
$ws_CUL.Range("I130").Value = 1000
$ws_CUL.Range("J130").Value = 2625
$ws_CUL.Range("K130").Value = 3000
$ws_CUL.Range("L130").Value = 7875
$ws_CUL.Range("M130").Value = 2020
$ws_CUL.Range("N130").Value = -17915
$ws_CUL.Range("H131").Value = 585.2759
$ws_CUL.Range("I131").Value = 518.92
$ws_CUL.Range("J131").Value = 1000
$ws_CUL.Range("K131").Value = 1556.76
$ws_CUL.Range("L131").Value = 3000
$ws_CUL.Range("M131").Value = 3483.24
$ws_CUL.Range("N131").Value = -13080
$ws_CUL.Range("H137").Value = 4862890
$ws_CUL.Range("I137").Value = 65772.3
$ws_CUL.Range("J137").Value = 25250642
$ws_CUL.Range("K137").Value = 197316.9
$ws_CUL.Range("L137").Value = 75751926
$ws_CUL.Range("M137").Value = -192216.9
$ws_CUL.Range("N137").Value = -75762126

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 1865.5
$ws_LTW.Range("I7").Value = 1477
$ws_LTW.Range("J7").Value = 2642.5
$ws_LTW.Range("K7").Value = 1477
$ws_LTW.Range("L7").Value = 2642.5
$ws_LTW.Range("M7").Value = -1365
$ws_LTW.Range("N7").Value = -2866.5
$ws_LTW.Range("H40").Value = 1960.909
$ws_LTW.Range("I40").Value = 1724.2858
$ws_LTW.Range("K40").Value = 1724.2858
$ws_LTW.Range("M40").Value = -1588.2858
$ws_LTW.Range("H61").Value = 1520.6538
$ws_LTW.Range("I61").Value = 1228
$ws_LTW.Range("K61").Value = 1228
$ws_LTW.Range("M61").Value = -1026
$ws_LTW.Range("H113").Value = 1520.6538
$ws_LTW.Range("I113").Value = 1228
$ws_LTW.Range("K113").Value = 1228
$ws_LTW.Range("M113").Value = 942
$ws_LTW.Range("H126").Value = 1865.5
$ws_LTW.Range("I126").Value = 1477
$ws_LTW.Range("J126").Value = 2642.5
$ws_LTW.Range("K126").Value = 4431
$ws_LTW.Range("L126").Value = 7927.5
$ws_LTW.Range("M126").Value = -1961
$ws_LTW.Range("N126").Value = -12867.5
$ws_LTW.Range("H132").Value = 5500.6943
$ws_LTW.Range("I132").Value = 9128.9375
$ws_LTW.Range("J132").Value = 2598.1
$ws_LTW.Range("K132").Value = 27386.8125
$ws_LTW.Range("L132").Value = 7794.299999999999
$ws_LTW.Range("M132").Value = -24856.8125
$ws_LTW.Range("N132").Value = -12854.3
$ws_LTW.Range("H136").Value = 2425.1792
$ws_LTW.Range("I136").Value = 888.9322
$ws_LTW.Range("J136").Value = 13755
$ws_LTW.Range("K136").Value = 2666.7966
$ws_LTW.Range("L136").Value = 41265
$ws_LTW.Range("M136").Value = -116.7965999999997
$ws_LTW.Range("N136").Value = -46365

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H107").Value = 650.2
$ws_WVR.Range("I107").Value = 587.75
$ws_WVR.Range("J107").Value = 900
$ws_WVR.Range("K107").Value = 1763.25
$ws_WVR.Range("L107").Value = 2700
$ws_WVR.Range("M107").Value = 156.75
$ws_WVR.Range("N107").Value = -6540
$ws_WVR.Range("H122").Value = 48982004
$ws_WVR.Range("J122").Value = 47621670
$ws_WVR.Range("L122").Value = 142865010
$ws_WVR.Range("N122").Value = -142869910
$ws_WVR.Range("H136").Value = 781.9091
$ws_WVR.Range("I136").Value = 376.8889
$ws_WVR.Range("J136").Value = 1490.6945
$ws_WVR.Range("K136").Value = 1130.6667
$ws_WVR.Range("L136").Value = 4472.083500000001
$ws_WVR.Range("M136").Value = 1419.3333
$ws_WVR.Range("N136").Value = -9572.083500000001
